# Fix header labels: remove embedded line breaks so the headers read as
# plain single-line text instead of wrapping onto two lines.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Parent Segment ID"
$ws.Range("D1").Value = "Segment Description"

# Leave the active selection on D1, matching the saved view state.
$ws.Range("D1").Select()
